# "Fill empty predicted prices"
# C2 previously held a hard-coded value (5386); replace it with the
# formula that actually computes "TimeTaken in Hours" from
# "TimeTaken in Minutes" (B2), i.e. B2 / 60.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Formula = "=B2/60"

# Leave the sheet with the same active selection Excel saved with
# after making this edit.
$ws.Range("C3").Select()
